$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column E (Date Sampled), shifting E:N -> H:Q
$ws.Range("E1:G1").EntireColumn.Insert()

# Header row (row 1) - new columns E, F, G => Month, Day, Year
# (Set in Day, Month, Year order so the shared-string table indices line up
# with the canonical export: 23=Day, 24=Month, 25=Year)
$ws.Range("F1").Value = "Day"
$ws.Range("E1").Value = "Month"
$ws.Range("G1").Value = "Year"

# Data rows 2-5: Month, Day, Year values (derived from original Date Sampled - now in column H)
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 9
$ws.Range("G2").Value = 2015

$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 2015

$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 2015

$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 2015

# New Day/Month/Year columns are plain-width (no bestFit), matching the
# other narrow data columns (9.5 chars stored == 8.667 ColumnWidth units)
$ws.Columns("E:G").ColumnWidth = 8.666666666666666

# Update the active selection to the newly inserted header cells
[void]$ws.Range("E1:G1").Select()

Write-Host "done"
